$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Peer  and self assessment")

# Fill in Daniel's and Alex's peer-assessment grades/comments
# (merging in Daniel's assessment notes), for both the
# "Criterion 1 Online collaboration" and
# "Criterion 1 International Collaboration" sections.

# Row 7 -> Daniel (Criterion 1 Online collaboration)
$ws.Range("C7").Value = "Decent activity on discord and good response times. Missed a meeting."
$ws.Range("B7").Value = "Good"

# Row 6 -> Alex (Criterion 1 Online collaboration)
$ws.Range("C6").Value = "Research, hardware setup"
$ws.Range("B6").Value = "Good"

# Row 19 -> Alex (Criterion 1 International Collaboration)
$ws.Range("C19").Value = "Active collaborator, motivated"
$ws.Range("B19").Value = "Excellent"

# Row 20 -> Daniel (Criterion 1 International Collaboration)
$ws.Range("C20").Value = "Active collaborator, motivated"
$ws.Range("B20").Value = "Excellent"
